$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet ---
# Title
$metadata.Range("B5").Value = "DMI Internal Diameter"
# Date
$metadata.Range("B8").Value = "2026-02-25T08:15:31+00:00"
# Description
$metadata.Range("B12").Value = "Extension créée dans ce volet pour représenter le diamètre interne."

# --- Elements sheet (root Extension row) ---
# Short
$elements.Range("L2").Value = "DMI Internal Diameter"
# Definition
$elements.Range("M2").Value = "Extension créée dans ce volet pour représenter le diamètre interne."
# Mapping: RIM Mapping
$elements.Range("AK2").Value = ""
